# Generate Report for Handoff
#
# Refreshes the localization status report to reflect a newly generated
# handoff run: a new source-file GUID, a new xliff content hash, and
# refreshed handoff/generation timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "4d1a104d-c357-45b1-8b9b-17f842cf4a20"
$newGuid = "3eea3430-8d03-49a3-b58e-ce8bc5516f1a"

$oldHash = "854193947a2c2777ff2f5c10d0d3955173a72d1b"
$newHash = "b137fb26de8049eb53bdb2f598052a9a5a8cb5e5"

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/df2b63ecd11c260a0b2488ae3735cf6c6aed5380/e2e/"
# The hyperlink target itself is unchanged by this commit; only the
# cell text / hyperlink display label is refreshed to the new file name.
$fileAddress = $repoBase + $oldGuid + ".md"

function Update-Hyperlink($ws, $cell, $address, $displayText) {
    $cell.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($cell, $address, "", "", $displayText) | Out-Null
}

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"

$overviewDisplay = "e2e\$newGuid.md"
Update-Hyperlink $wsOverview $wsOverview.Range("B2") $fileAddress $overviewDisplay

$wsOverview.Range("G2").Value = "2016-08-29 09:01:57"

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhCnDisplay = "$newGuid.md"
Update-Hyperlink $wsZhCn $wsZhCn.Range("A2") $fileAddress $zhCnDisplay

$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-29 09:01:53"

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")

$deDeDisplay = "$newGuid.md"
Update-Hyperlink $wsDeDe $wsDeDe.Range("A2") $fileAddress $deDeDisplay

$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-29 09:01:57"
